$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Sheet1"

# Remove existing merged regions and clear old content (keep base styling)
$ws.Cells.UnMerge()
$ws.Range("A1:L11").ClearContents()

# Resize columns 1 and 2 (closest achievable values through the ColumnWidth API)
$ws.Columns.Item(1).ColumnWidth = 26
$ws.Columns.Item(2).ColumnWidth = 21.666666666666668

# Extend the bordered/wrapped formatting down to the new rows (12-15)
$ws.Range("A1").Copy()
$ws.Range("A12:L15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Header / title block ----
$ws.Range("A1").Value = "MASTER PACKAGE"
$ws.Range("A2").Value = "WesternGlove Centric8 PROD"
$ws.Range("B2").Value = "M12225BVS563:KONRAD"
$ws.Range("C2").Value = "BOM"
$ws.Range("D2").Value = "MASTER"
$ws.Range("A3").Value = "Placements"

# ---- Table header row ----
$ws.Range("A5").Value = "Code"
$ws.Range("B5").Value = "Product"
$ws.Range("C5").Value = "Type"
$ws.Range("D5").Value = "Description :"
$ws.Range("E5").Value = "Main`nMaterial"
$ws.Range("F5").Value = "Composition"
$ws.Range("G5").Value = "Coating`nComposition"
$ws.Range("H5").Value = "DUNE WASH"
$ws.Range("I5").Value = "Weight`n/ Yield"
$ws.Range("J5").Value = "Common`nQty"
$ws.Range("K5").Value = "Image"
$ws.Range("L5").Value = "Supplier"

# ---- WASH (1) group ----
$ws.Range("A6").Value = "WASH (1)"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "563"
$ws.Range("B7").Value = "BVS563"
$ws.Range("C7").Value = "Wash"
$ws.Range("D7").Value = "M12225BVS563"

# ---- POCKET SCROLL (1) group ----
$ws.Range("A8").Value = "POCKET SCROLL (1)"

$ws.Range("A9").Value = "ICON"
$ws.Range("B9").Value = "SILVER ICON`nEMBROIDERY"
$ws.Range("C9").Value = "Trim"
$ws.Range("D9").Value = "TINY LOCKSTITCH`n(TEX 27x1, 20 SPI)`n4 ROW STITCHING`n@ VERTICAL`nEMBROIDERY`n4 ROW STITCHING`n@ DIAGONAL`nEMBROIDERY"
$ws.Range("H9").Value = "ICON GUNMETAL`nW32633 (A&E):`n500S"
$ws.Range("K9").Value = "3/8"
$ws.Range("L9").Value = "SILVERMOON`nJEANS`nLIMITED"

# ---- FABRIC (1) group ----
$ws.Range("A10").Value = "FABRIC (1)"

$ws.Range("A11").Value = "A831D9-22`nPOWER`nSTRETCH"
$ws.Range("B11").Value = "BVS"
$ws.Range("C11").Value = "Fabric"
$ws.Range("D11").Value = "BVS VIETNAM`nVERSION, testing to`nconfirm okay to use"
$ws.Range("F11").Value = "57% Cotton,`n22%`nPolyester,`n19% Lyocell,`n2% Elastane"
$ws.Range("H11").Value = "black/black"
$ws.Range("I11").Value = "9.7 oz"

$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = "1.9"
$ws.Range("L11").Value = "ADVANCE`nVIETNAM"

# ---- INSIDE TRIMS (5) group ----
$ws.Range("A12").Value = "INSIDE TRIMS (5)"

$ws.Range("A13").Value = "MEN'S FUSING"
$ws.Range("B13").Value = "MEN'S`nFUSING"
$ws.Range("C13").Value = "Fabric"
$ws.Range("D13").Value = "FUSE INSIDE TOP`nAND BOTTOM`nWAISTBAND (2 1/2""`nLONG) FOR SHANK`nREINFORCEMENT"

$ws.Range("A14").Value = "P3026"
$ws.Range("B14").Value = "ASMARA`nVIETNAM`nPOCKETING -`nP3026"
$ws.Range("C14").Value = "Fabric"
$ws.Range("F14").Value = "65%`nPolyester,`n35% Cotton"
$ws.Range("H14").Value = "BLACK"
$ws.Range("I14").Value = "95GSM"
$ws.Range("L14").Value = "Asmara`nVietnam"

$ws.Range("A15").Value = "Displaying 1 - 5 of 20 results"

# Re-create the merged banner rows
$ws.Range("A6:L6").Merge()
$ws.Range("A8:L8").Merge()
$ws.Range("A10:L10").Merge()
$ws.Range("A12:L12").Merge()
$ws.Range("A15:L15").Merge()

# Normalize formatting across the whole used range so every written cell
# shares the same (bordered, wrapped, left/top-aligned) style, and so any
# transient number-format tweaks above are cleared back to the common style.
$ws.Range("A1").Copy()
$ws.Range("A1:L15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Select()
